# Combined LHS_Gen.rb and Morris.rb -- Meters sheet clean-up.
$wb = $excel.ActiveWorkbook

$wsMeters = $wb.Worksheets.Item("Meters")
$wsSources = $wb.Worksheets.Item("Sources")

# Update the meter / report-frequency pairs that survive.
$wsMeters.Range("B2").Value = "Monthly"
$wsMeters.Range("A3").Value = "Gas:Facility"
$wsMeters.Range("B3").Value = "Monthly"

# Remove the two now-unneeded meter rows (Carbon Equivalent related entries).
$wsMeters.Rows.Item(4).Delete()
$wsMeters.Rows.Item(4).Delete()

# Leave the selection where Excel would land after deleting rows 4 and 5.
$wsMeters.Range("A6").Select()

# Make Sources the active sheet/tab (was TotalEnergy before).
$wsSources.Activate()
